$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.914.72'
$ws.Range("E2").Value = '  -3.36%  '
$ws.Range("D3").Value = '3.229.47'
$ws.Range("E3").Value = '  -3.79%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("E5").Value = '  -5.66%  '
$ws.Range("E6").Value = '  -10.01%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = '3.230.63'
$ws.Range("E8").Value = '  -3.81%  '
$ws.Range("E9").Value = '  -4.47%  '
$ws.Range("E10").Value = '  -4.66%  '
$ws.Range("E11").Value = '  -5.69%  '
$ws.Range("E12").Value = '  -5.15%  '
$ws.Range("D13").Value = '3.781.88'
$ws.Range("E13").Value = '  -3.94%  '
$ws.Range("E14").Value = '  -1.17%  '
$ws.Range("E15").Value = '  -7.81%  '
$ws.Range("D16").Value = '3.228.06'
$ws.Range("E16").Value = '  -5.27%  '
$ws.Range("E17").Value = '  -6.25%  '
$ws.Range("D18").Value = '59.039.77'
$ws.Range("E18").Value = '  -3.32%  '
$ws.Range("E19").Value = '  -7.03%  '
$ws.Range("E20").Value = '  -6.73%  '
$ws.Range("E21").Value = '  -6.57%  '
$ws.Range("E22").Value = '  -3.34%  '
$ws.Range("E23").Value = '  -0.10%  '
$ws.Range("E24").Value = '  -5.77%  '
$ws.Range("E25").Value = '  -7.58%  '
$ws.Range("D26").Value = '3.361.82'
$ws.Range("E26").Value = '  -4.81%  '
$ws.Range("E27").Value = '  -2.73%  '
$ws.Range("D28").Value = '0.0₃0974'
$ws.Range("E28").Value = '  -10.05%  '
$ws.Range("E29").Value = '  -0.28%  '
$ws.Range("E30").Value = '  -4.70%  '
$ws.Range("E31").Value = '  -0.09%  '
$ws.Range("E32").Value = '  -7.34%  '
$ws.Range("E33").Value = '  -7.94%  '
$ws.Range("E34").Value = '  -3.77%  '
$ws.Range("E35").Value = '  -3.29%  '
$ws.Range("E36").Value = '  -3.56%  '
$ws.Range("E37").Value = '  -8.56%  '
$ws.Range("E38").Value = '  -5.85%  '
$ws.Range("E39").Value = '  -6.46%  '
$ws.Range("E40").Value = '  -13.52%  '
$ws.Range("E41").Value = '  -6.35%  '
$ws.Range("D42").Value = '3.260.46'
$ws.Range("E42").Value = '  -3.93%  '
$ws.Range("E43").Value = '  -3.01%  '
$ws.Range("E44").Value = '  -5.91%  '
$ws.Range("E45").Value = '  -6.68%  '
$ws.Range("E46").Value = '  -4.56%  '
$ws.Range("E47").Value = '  -6.78%  '
$ws.Range("E48").Value = '  -0.04%  '
$ws.Range("D49").Value = '2.287.38'
$ws.Range("E49").Value = '  -8.95%  '
$ws.Range("E50").Value = '  -6.38%  '
$ws.Range("E51").Value = '  -9.45%  '

# Force the following price values to remain plain text (not auto-converted to numbers)
# by building them via a helper cell formula and pasting the computed value back as text.
$ws.Range("Z1").Formula = '="536.74"'
$ws.Range("Z1").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="135.77"'
$ws.Range("Z1").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="7.58"'
$ws.Range("Z1").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="0.115"'
$ws.Range("Z1").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="0.395"'
$ws.Range("Z1").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="26.03"'
$ws.Range("Z1").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="8.28"'
$ws.Range("Z1").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="361.13"'
$ws.Range("Z1").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="70.71"'
$ws.Range("Z1").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="7.07"'
$ws.Range("Z1").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="1.93"'
$ws.Range("Z1").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="7.09"'
$ws.Range("Z1").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="22.04"'
$ws.Range("Z1").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="163.69"'
$ws.Range("Z1").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="4.93"'
$ws.Range("Z1").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="25.97"'
$ws.Range("Z1").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="41.04"'
$ws.Range("Z1").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="0.717"'
$ws.Range("Z1").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="4.01"'
$ws.Range("Z1").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="6.28"'
$ws.Range("Z1").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="20.75"'
$ws.Range("Z1").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$excel.CutCopyMode = $false
